$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 41

$ws.Range("E17").Value = 94

$ws.Range("E18").Value = 99
$ws.Range("F18").Value = 39
$ws.Range("H18").Value = 39

$ws.Range("E19").Value = 43

$ws.Range("E26").Value = 23
$ws.Range("F26").Value = 10
$ws.Range("H26").Value = 10

$ws.Range("E36").Value = 83

$ws.Range("E38").Value = 60

$ws.Range("E39").Value = 21
$ws.Range("F39").Value = 12
$ws.Range("H39").Value = 12

$ws.Range("E40").Value = 17
$ws.Range("F40").Value = 9
$ws.Range("H40").Value = 9

$ws.Range("E43").Value = 20

$ws.Range("E49").Value = 54

$ws.Range("E52").Value = 2

$ws.Range("E57").Value = 11

$ws.Range("E67").Value = 35
$ws.Range("F67").Value = 19
$ws.Range("H67").Value = 19

$ws.Range("E76").Value = 42

$ws.Range("E79").Value = 28
$ws.Range("F79").Value = 11
$ws.Range("H79").Value = 11

$ws.Range("E82").Value = 11

$ws.Range("E89").Value = 28
